$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the "2020"-prefixed parameter labels to "2022" (model now runs with
# 2022 technology costs). Each of these cells holds rich text: a plain run
# with the year, followed by a separately-formatted run with the rest of the
# label (e.g. " CapEx"). We only touch the year run's text, then re-assert
# the trailing run's font so the writer keeps it as a distinct formatted run
# instead of flattening the cell to plain text.
function Set-YearLabel($cellAddress) {
    $cell = $ws.Range($cellAddress)
    $fullText = $cell.Text
    $suffix = $fullText.Substring(4)
    $cell.Characters(1, 4).Text = "2022"
    $rest = $cell.Characters(5, $suffix.Length)
    $rest.Font.Name = "Calibri"
    $rest.Font.Size = 11
    $rest.Font.Bold = $false
}

Set-YearLabel "A9"
Set-YearLabel "A15"
Set-YearLabel "A22"
Set-YearLabel "A28"

# Restore the current selection/scroll position captured in the edit.
$ws.Range("A28").Select()
